# Bài 2 - Phần 2
#
# The only content change in this commit is on the "4. Vẽ hình đa giác
# đều" slide (slide 5): the first bullet of the body placeholder gets an
# extra sentence appended ("Các cạnh này có kích thước bằng nhau."),
# while the other three bullets stay exactly as they were.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$oldFirstLine = "Là hình có n cạnh."
$newFirstLine = "Là hình có n cạnh. Các cạnh này có kích thước bằng nhau."

# Read each paragraph's text individually first (TextRange.Paragraphs(i,1).Text
# includes the trailing paragraph mark, except for the very last paragraph),
# so the other bullets are carried over verbatim.
$paraCount = $tr.Paragraphs().Count
$lines = @()
for ($i = 1; $i -le $paraCount; $i++) {
    $lines += $tr.Paragraphs($i, 1).Text
}

if ($lines.Count -gt 0 -and $lines[0] -eq ($oldFirstLine + "`r")) {
    $lines[0] = $newFirstLine + "`r"
} elseif ($lines.Count -gt 0 -and $lines[0] -eq $oldFirstLine) {
    $lines[0] = $newFirstLine
}

$newFull = [string]::Join("", $lines)

# Clear the whole placeholder before re-assigning the full text: editing
# TextRange.Text in place does a run-level diff against the previous
# contents and can split the edited paragraph into multiple <a:r> runs.
# Clearing first makes every paragraph come back as a single plain run,
# matching the original (and target) formatting.
$tr.Text = ""
$tr.Text = $newFull
